$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# --- Sheet2: simple login (userName / password) form, with the password value
#     typed as an email-like string that Excel auto-hyperlinks ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1:B2").NumberFormat = "@"
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "faruq"
$ws2.Range("B2").Value = "faruq@123"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:faruq@123")

# --- Sheet3: checking / saving account transfer amounts ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "checking "
$ws3.Range("B1").Value = "saving"
$ws3.Range("A2").Value = 500
$ws3.Range("B2").Value = 5000
$ws3.Range("A3").Value = 300
$ws3.Range("B3").Value = 4000

$null = $ws3.Activate()
$null = $ws3.Rows(2).Select()

# --- Sheet1: update LastName sample data + add a numeric cell, then fit two new columns ---
$ws1.Range("B2").Value = "Molla"
$ws1.Range("B3").Value = "Fima"
$ws1.Range("B4").Value = "Nisa"
$ws1.Range("I2").Value = 234

$ws1.Columns.Item(5).ColumnWidth = 8.1
$ws1.Columns.Item(6).ColumnWidth = 6

# Leave Sheet1 as the active / selected tab, matching the saved workbook state.
$null = $ws1.Activate()
$null = $ws1.Range("B4").Select()
